# Add a "script value" example: a new managed-text line that prints the
# values of two custom script variables (t_Foo / t_Bar), plus a new
# "Text>Script" sheet that declares those variables' example values.

$wb = $excel.ActiveWorkbook

# --- 1. Scripts>Script001: insert the new line before "@goto Script002" ---
$ws1 = $wb.Worksheets.Item(1)

# Row 16 currently holds "\n@goto Script002" in column A only. Push that
# down to a brand-new row 17 (column A only, default style) ...
$ws1.Range("A17").Value = "`n@goto Script002"
$ws1.Range("A17").Style = "Normal"
$ws1.Rows.Item(17).AutoFit()

# ... and replace row 16 with the new example line + its translations.
$ws1.Range("A16").Value = "`n{0}`n"
$ws1.Range("B16").Value = "Value of t_Foo: {t_Foo}, value of t_Bar: {t_Bar}."
$ws1.Range("C16").Value = " t_Foo の値: {t_Foo},  t_Bar の値: {t_Bar}."
$ws1.Range("D16").Value = "Значение t_Foo: {t_Foo}, значениеf t_Bar: {t_Bar}."
$ws1.Range("B16:D16").Style = "Normal"
$ws1.Rows.Item(16).AutoFit()

# --- 2. New "Text>Script" sheet declaring the t_Foo / t_Bar example values ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$scriptSheet = $wb.Worksheets.Add($null, $lastSheet)
$scriptSheet.Name = "Text>Script"

$scriptSheet.Range("A1").Value = "Template"
$scriptSheet.Range("B1").Value = "Arguments"
$scriptSheet.Range("C1").Value = "ja"
$scriptSheet.Range("D1").Value = "ru"

$scriptSheet.Range("A2").Value = "t_Foo: {0}"
$scriptSheet.Range("B2").Value = "Foo"
$scriptSheet.Range("C2").Value = "ふぉお"
$scriptSheet.Range("D2").Value = "Фу"

$scriptSheet.Range("A3").Value = "t_Bar: {0}"
$scriptSheet.Range("B3").Value = "Bar"
$scriptSheet.Range("C3").Value = "ばら"
$scriptSheet.Range("D3").Value = "Бар"

$scriptSheet.Range("A18").Select()

# --- 3. Restore Script001 as the active/selected sheet ---
$ws1.Select()
$ws1.Range("C18").Select()
